# Saldo_guide.xlsx - "Add files via upload" re-upload edit
# The workbook is a daily client-balance export. This re-upload:
#   1. Renames the (only) worksheet to the next day's export tag
#      (IClientBalance-20241106-112141- -> IClientBalance-20241107-101553-)
#   2. Rolls the "Dt. Referencia" column (G, rows 2-274) from 2024-11-06
#      (serial 45602) to 2024-11-07 (serial 45603)
#   3. Corrects a single data-entry: row 264 (account 58420177172) had its
#      Saldo Previsto / Vl. Total (E264 / H264) mis-keyed as 964.68;
#      corrected to 164.68

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet to match the new export timestamp.
$ws.Name = "IClientBalance-20241107-101553-"

# 2) Bump every "Dt. Referencia" cell in the data body one day forward.
$ws.Range("G2:G274").Value = 45603

# 3) Fix the mis-keyed balance on row 264.
$ws.Range("E264").Value = 164.68
$ws.Range("H264").Value = 164.68
